$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-12-23"

# Update the header label in I1 ("2022 (through 12-22)" -> "2022 (through 12-23)")
$ws.Range("I1").Value = "2022 (through 12-23)"

# Update December (row 13) and Total (row 14) values in column I
$ws.Range("I13").Value = 99
$ws.Range("I14").Value = 1616
